$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(12.40453383752919, 8.440546529063088, 14.0160008166415, 14.7726355906313, 0, 42.35544593076597, 17.85648791140229, 28.75904256296874, 8.88824659559746, 9.523861557817211, 12.11450060032691, 15.48203246089518, 0, 28.95271429876113),
    @(12.19519397302061, 8.413745393550723, 14.02069832662878, 14.80312046728939, 0, 42.49481110824836, 17.90885063336371, 28.85983927469492, 8.895390369357141, 9.36666126740276, 12.12484038866509, 15.45352418494169, 0, 29.0458321197443),
    @(12.06653563793971, 8.39719269210102, 14.0259115380325, 14.82348445052418, 0, 42.58955865326764, 17.94329354353553, 28.92599147449038, 8.900012537860388, 9.269902490460321, 12.13253603460777, 15.43778123833013, 0, 29.10768577326957),
    @(12.01414013588727, 8.390424796961643, 14.0286227085006, 14.83219742720272, 0, 42.63047237691628, 17.95790621432549, 28.954021604072, 8.901955587558856, 9.23045947165679, 12.13601147115306, 15.43181320103583, 0, 29.13406806152121),
    @(12.00544385223971, 8.389299683735448, 14.02910836298021, 14.83366926046104, 0, 42.63740505739737, 17.96036749760528, 28.9587407874501, 8.902281826810642, 9.223910575958865, 12.13660908294047, 15.43084935661922, 0, 29.13851986942223),
    @(12.06582879085362, 8.397101506804397, 14.025945724782, 14.82360027779582, 0, 42.5901011097087, 17.9434882785027, 28.9263651548753, 8.900038501456573, 9.269370538825882, 12.13258153046434, 15.4376989343203, 0, 29.10803681092997),
    @(12.33241745701951, 8.431326635872583, 14.01713787772174, 14.78280545216714, 0, 42.40159191335152, 17.87406721698645, 28.79291298648094, 8.890660917414998, 9.469735856151932, 12.11778662033461, 15.4718400074741, 0, 28.98385013515412),
    @(12.85143735936938, 8.497597601187881, 14.01828861053364, 14.71584379888772, 0, 42.1049223231423, 17.75609311098123, 28.56500372078056, 8.87413524145852, 9.858764941896842, 12.09943030964053, 15.55256038142518, 0, 28.77745128040015),
    @(13.22685030010545, 8.545670180234589, 14.03027961435449, 14.67456118750162, 0, 41.93169940510114, 17.68045338254636, 28.41811377652819, 8.863119175582151, 10.13961085565686, 12.09239725203823, 15.61997260332394, 0, 28.6484529772109),
    @(13.39562585210626, 8.567383836233315, 14.0381322596869, 14.65749188895356, 0, 41.862657711405, 17.64843192076152, 28.35574412302843, 8.858349738214022, 10.26577227310452, 12.09058893905736, 15.6523343006414, 0, 28.59468646637175),
    @(13.45919168454928, 8.575582039933909, 14.04144836140598, 14.6512735494989, 0, 41.83792014682395, 17.63664901499144, 28.33276589695207, 8.856578277641901, 10.31327544322031, 12.09010326323683, 15.66482661600499, 0, 28.5750335139325),
    @(13.44551791802762, 8.573817522605543, 14.04071898746599, 14.65260187275198, 0, 41.84318518857165, 17.63917142881571, 28.33768621772552, 8.856958255968571, 10.30305750145804, 12.09019902301545, 15.66212570384998, 0, 28.57923466683995),
    @(13.40086270389891, 8.568058805767183, 14.03839822783092, 14.65697538655755, 0, 41.86059431740456, 17.64745566341057, 28.35384086938281, 8.858203306081242, 10.26968605952577, 12.09054499766101, 15.65335732439309, 0, 28.59305542928166),
    @(13.37346338918817, 8.564528202282753, 14.0370212227909, 14.65968623726977, 0, 41.87144125586467, 17.65257464428321, 28.36381937381526, 8.858970439511186, 10.24920853409047, 12.09078281613802, 15.64801719981886, 0, 28.60161316347878),
    @(13.21577512324567, 8.54424785377487, 14.02981448231413, 14.67571107944206, 0, 41.9364081345902, 17.68259406867665, 28.4222793269492, 8.863435722627349, 10.1313301601846, 12.0925433420913, 15.61789128638146, 0, 28.65206569313817),
    @(13.11848368138755, 8.53176543021568, 14.02600594298672, 14.68597950364273, 0, 41.97876532863286, 17.70162121642235, 28.45928249269629, 8.866236857756483, 10.05857625306896, 12.09397899951565, 15.59983996326593, 0, 28.68427595812891),
    @(13.06233910666639, 8.524571421653578, 14.02404122549905, 14.69204664903489, 0, 42.00404635167047, 17.7127898671071, 28.48098472962984, 8.867870767454287, 10.0165822112281, 12.09493574435836, 15.58961729610327, 0, 28.70326509309015),
    @(13.04329953711177, 8.522133233952482, 14.02341486401545, 14.69412855100166, 0, 42.01276366791601, 17.71660999119186, 28.48840470604344, 8.86842789678902, 10.00233965097647, 12.09528220761955, 15.58618373710346, 0, 28.70977393415939),
    @(13.12886010847128, 8.533095715544055, 14.02638800972222, 14.68486975116337, 0, 41.97416126554135, 17.69957248725919, 28.45530008029992, 8.865936316550938, 10.06633664880218, 12.09381262072708, 15.60174504446064, 0, 28.68079923557246),
    @(13.4139888558492, 8.569750955825629, 14.03907061623624, 14.6556841226986, 0, 41.85544261907636, 17.64501307897036, 28.34907849454495, 8.857836666264072, 10.2794957504941, 12.09043798093978, 15.65592641075109, 0, 28.58897673826822),
    @(13.59829647922916, 8.593564766936842, 14.04935440004626, 14.63803995544735, 0, 41.7860559113091, 17.61135401999466, 28.28338572395482, 8.852744809453943, 10.41720753691316, 12.08939245269466, 15.6927192093735, 0, 28.53308805504378),
    @(13.50013315429366, 8.580868599128964, 14.04368404095187, 14.64732627718984, 0, 41.82233716466351, 17.62913573928013, 28.31810606016844, 8.85544402070097, 10.34386785595973, 12.08984465518963, 15.67295783159726, 0, 28.56253954154153),
    @(13.12416957678838, 8.532494348521922, 14.02621457676425, 14.68537096034323, 0, 41.97623986802178, 17.70049800221063, 28.45709919378112, 8.866072118019339, 10.06282869039702, 12.09388743131794, 15.60088327311157, 0, 28.68236959465605),
    @(12.7118123722676, 8.479770401052173, 14.01601298065036, 14.73256644649621, 0, 42.17733890808918, 17.78606771216133, 28.62304644012665, 8.878407474070706, 9.754212631973971, 12.10325959320656, 15.52927789563301, 0, 28.82931177496449),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $excelCol = $j + 2
        $ws.Cells.Item($excelRow, $excelCol).Value = $row[$j]
    }
}
